$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths (29.14 chars rounds/saves to an exact stored width of 30)
$ws.Columns.Item(1).ColumnWidth = 29.14
$ws.Columns.Item(2).ColumnWidth = 29.14
$ws.Columns.Item(3).ColumnWidth = 29.14
$ws.Columns.Item(4).ColumnWidth = 29.14

# Header row
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "peso"
$ws.Range("C1").Value = "zone"
$ws.Range("D1").Value = "partenza"

# Row 2
$ws.Range("A2").Value = "40279 (interno)"
$ws.Range("B2").Value = "CAMPO VUOTO"
$ws.Range("C2").Value = "CAMPO VUOTO"
$ws.Range("D2").Value = 0

# Row 3
$ws.Range("A3").Value = "40176 (interno)"
$ws.Range("B3").Value = "CAMPO VUOTO"
$ws.Range("C3").Value = "CAMPO VUOTO"
$ws.Range("D3").Value = 0

# Highlight the "campo vuoto" (empty field) cells with a solid yellow fill
$ws.Range("B2:C3").Interior.Color = 65535
